$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 322 (pushes the existing rows 322..343 down to 324..345)
$ws.Rows.Item(322).Insert()
$ws.Rows.Item(322).Insert()

# --- New row 322 ("Primera") ---
$ws.Cells.Item(322, 1).Value = 7
$ws.Cells.Item(322, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(322, 3).Value = "Ñuble"
$ws.Cells.Item(322, 4).Value = 44931
$ws.Cells.Item(322, 5).Value = 16
$ws.Cells.Item(322, 6).Value = 100112009
$ws.Cells.Item(322, 7).Value = "Acelga"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 400
$ws.Cells.Item(322, 11).Value = 600
$ws.Cells.Item(322, 12).Value = 700
$ws.Cells.Item(322, 13).Value = 650
$ws.Cells.Item(322, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(322, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(322, 16).Value = 650
$ws.Cells.Item(322, 17).Value = 1
$ws.Cells.Item(322, 18).Value = "Hortaliza"

# --- New row 323 ("Segunda") ---
$ws.Cells.Item(323, 1).Value = 7
$ws.Cells.Item(323, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(323, 3).Value = "Ñuble"
$ws.Cells.Item(323, 4).Value = 44931
$ws.Cells.Item(323, 5).Value = 16
$ws.Cells.Item(323, 6).Value = 100112009
$ws.Cells.Item(323, 7).Value = "Acelga"
$ws.Cells.Item(323, 8).Value = "Sin especificar"
$ws.Cells.Item(323, 9).Value = "Segunda"
$ws.Cells.Item(323, 10).Value = 300
$ws.Cells.Item(323, 11).Value = 500
$ws.Cells.Item(323, 12).Value = 500
$ws.Cells.Item(323, 13).Value = 500
$ws.Cells.Item(323, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(323, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(323, 16).Value = 500
$ws.Cells.Item(323, 17).Value = 1
$ws.Cells.Item(323, 18).Value = "Hortaliza"
